# "Xong automation Them thuoc cho Benh nhan"
# Adds a "Result" (T = Thanh cong / pass-marker) column for the medication
# rows, and fills in the Pass/Fail counters as text values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New D column values ("T") for rows 2-7
$ws.Range("D2").Value = "T"
$ws.Range("D3").Value = "T"
$ws.Range("D4").Value = "T"
$ws.Range("D5").Value = "T"
$ws.Range("D6").Value = "T"
$ws.Range("D7").Value = "T"

# Pass / Fail counts stored as text ("6" / "0"), matching the style already
# used by column B (reuse existing General-formatted style rather than
# switching the cell to a Text number format).
$ws.Range("B10").Formula = "=TEXT(6,""0"")"
$ws.Range("B10").Copy()
$ws.Range("B10").PasteSpecial(-4163)

$ws.Range("B11").Formula = "=TEXT(0,""0"")"
$ws.Range("B11").Copy()
$ws.Range("B11").PasteSpecial(-4163)

$excel.CutCopyMode = $false

# Best-effort: mark the workbook for full recalculation / full precision on
# load (matches calcPr@fullCalcOnLoad / calcPr@fullPrecision in the target).
$wb.PrecisionAsDisplayed = $true
$wb.ForceFullCalculation = $true
